# Trade #41 closed at 2026-02-17 08:33:29 - unknown UNKNOWN +0.000%

$wb = $excel.ActiveWorkbook

# --- Summary sheet -------------------------------------------------
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B3").Value = 1199.78   # Current Capital
$wsSummary.Range("B4").Value = -0.22     # Total P&L $
$wsSummary.Range("B6").Value = 41        # Total Trades
$wsSummary.Range("B7").Value = 15        # Winning Trades
$wsSummary.Range("B9").Value = 36.59     # Win Rate %

# --- Strategy Status sheet ------------------------------------------
$wsStatus = $wb.Worksheets.Item("Strategy Status")
$wsStatus.Range("C4").Value = 99.78      # MarketMaking Capital
$wsStatus.Range("D4").Value = 41         # MarketMaking Trades
$wsStatus.Range("E4").Value = -0.22      # MarketMaking P&L $
$wsStatus.Range("F4").Value = -0.22      # MarketMaking P&L %
$wsStatus.Range("G4").Value = 36.59      # MarketMaking Win Rate %

# --- New trade row (#41) appended to "All Trades" and "MarketMaking" sheets
$newRow = @{
    A = 41
    B = "2026-02-17"
    C = "08:33:23"
    D = "MarketMaking"
    E = "UP"
    F = 0.44
    G = 0.45
    H = "CLOSED"
    I = 2.2727
    J = 0.01
    K = 99.78
    L = 0
    M = 0
    N = 0.6
    O = "Normal spread capture: 19600 bps"
    P = "early_exit"
    Q = 0.13
}

foreach ($sheetName in @("All Trades", "MarketMaking")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("A42").Value = $newRow.A

    # Column B holds a date-shaped string ("2026-02-17") that must stay a
    # plain text value (matching the rest of the column, which is
    # inlineStr text, not a real date serial). Force text interpretation,
    # then drop the now-unneeded explicit "Text" number format so the
    # cell keeps using the sheet's default (un-styled) cell.
    $ws.Range("B42").NumberFormat = "@"
    $ws.Range("B42").Value = $newRow.B
    $ws.Range("B42").ClearFormats()

    $ws.Range("C42").Value = $newRow.C
    $ws.Range("D42").Value = $newRow.D
    $ws.Range("E42").Value = $newRow.E
    $ws.Range("F42").Value = $newRow.F
    $ws.Range("G42").Value = $newRow.G
    $ws.Range("H42").Value = $newRow.H
    $ws.Range("I42").Value = $newRow.I
    $ws.Range("J42").Value = $newRow.J
    $ws.Range("K42").Value = $newRow.K
    $ws.Range("L42").Value = $newRow.L
    $ws.Range("M42").Value = $newRow.M
    $ws.Range("N42").Value = $newRow.N
    $ws.Range("O42").Value = $newRow.O
    $ws.Range("P42").Value = $newRow.P
    $ws.Range("Q42").Value = $newRow.Q
}
